# Updates cryptocurrency price/volume data cells per the scraped refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds values that look numeric ("47.558.22", "10.00", ...)
# but must stay plain text, matching the source data. Force text storage via
# NumberFormat="@" before assignment, then drop back to the default "Normal"
# style afterwards so no stray formatting is left on the cells.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D28", "D29", "D30", "D33", "D35", "D38", "D40", "D42", "D43", "D44", "D45", "D48", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$priceValues = @{
    "D2" = '47.558.22'
    "D3" = '2.488.87'
    "D5" = '322.57'
    "D6" = '105.32'
    "D7" = '0.525'
    "D10" = '37.89'
    "D11" = '0.0815'
    "D12" = '0.124'
    "D13" = '18.27'
    "D14" = '7.14'
    "D15" = '2.876.81'
    "D16" = '2.501.77'
    "D17" = '0.843'
    "D18" = '47.410.72'
    "D19" = '12.67'
    "D20" = '6.56'
    "D21" = '0.0₃0936'
    "D22" = '70.67'
    "D23" = '250.79'
    "D24" = '2.38'
    "D26" = '26.14'
    "D28" = '2.29'
    "D29" = '10.00'
    "D30" = '35.01'
    "D33" = '19.85'
    "D35" = '0.0780'
    "D38" = '4.62'
    "D40" = '2.26'
    "D42" = '121.23'
    "D43" = '21.01'
    "D44" = '0.0297'
    "D45" = '1.962.77'
    "D48" = '9.20'
    "D51" = '79.53'
}
foreach ($addr in $priceCells) {
    $ws.Range($addr).Value = $priceValues[$addr]
}
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining changed cells (coin name, link, volume%) are plain text already
# and can be assigned directly.
$ws.Range("E2").Value = '  +5.01%  '
$ws.Range("E3").Value = '  +2.64%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +1.47%  '
$ws.Range("E6").Value = '  +2.89%  '
$ws.Range("E7").Value = '  +1.61%  '
$ws.Range("E9").Value = '  +2.35%  '
$ws.Range("E10").Value = '  +6.45%  '
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("E15").Value = '  +2.51%  '
$ws.Range("E16").Value = '  +2.34%  '
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("E18").Value = '  +4.92%  '
$ws.Range("E19").Value = '  +3.72%  '
$ws.Range("E20").Value = '  +3.08%  '
$ws.Range("E21").Value = '  +1.77%  '
$ws.Range("E22").Value = '  +2.78%  '
$ws.Range("E23").Value = '  +2.77%  '
$ws.Range("E24").Value = '  +5.58%  '
$ws.Range("E25").Value = '  +2.35%  '
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("E28").Value = '  +4.54%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("E29").Value = '  +4.38%  '
$ws.Range("E30").Value = '  +6.46%  '
$ws.Range("E31").Value = '  +5.97%  '
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("E33").Value = '  -1.94%  '
$ws.Range("E34").Value = '  +2.62%  '
$ws.Range("E35").Value = '  +2.27%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  +3.52%  '
$ws.Range("E38").Value = '  +3.95%  '
$ws.Range("E39").Value = '  +4.51%  '
$ws.Range("E40").Value = '  +1.87%  '
$ws.Range("E41").Value = '  +1.63%  '
$ws.Range("E42").Value = '  -3.64%  '
$ws.Range("E43").Value = '  +1.48%  '
$ws.Range("E44").Value = '  +2.51%  '
$ws.Range("E45").Value = '  +1.46%  '
$ws.Range("E46").Value = '  +1.45%  '
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("E49").Value = '  -1.38%  '
$ws.Range("E50").Value = '  +12.04%  '
$ws.Range("E51").Value = '  +4.09%  '
